$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.291738
$ws.Range("H2").Value = 18.875214
$ws.Range("I2").Value = 0.5742845621220376
$ws.Range("J2").Value = 0.5742845621220376
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.7065936666666666
$ws.Range("N2").Value = 2.119781
$ws.Range("O2").Value = 0.005187843618793344
$ws.Range("P2").Value = 0.005187843618793344
$ws.Range("Q2").Value = 4.445702223125999
$ws.Range("R2").Value = 40.01132000813399
$ws.Range("S2").Value = 0.002979298500976342
$ws.Range("T2").Value = 0.002979298500976342

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.291738
$ws.Range("H3").Value = 18.875214
$ws.Range("I3").Value = 0.5742845621220376
$ws.Range("J3").Value = 0.5742845621220376
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 111.9320066666667
$ws.Range("N3").Value = 335.79602
$ws.Range("O3").Value = 0.8218100075305903
$ws.Range("P3").Value = 0.8218100075305903
$ws.Range("Q3").Value = 704.24685976092
$ws.Range("R3").Value = 6338.22173784828
$ws.Range("S3").Value = 0.4719528003222135
$ws.Range("T3").Value = 0.4719528003222135

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.291738
$ws.Range("H4").Value = 18.875214
$ws.Range("I4").Value = 0.5742845621220376
$ws.Range("J4").Value = 0.5742845621220376
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 23.563205
$ws.Range("N4").Value = 70.689615
$ws.Range("O4").Value = 0.1730021488506163
$ws.Range("P4").Value = 0.1730021488506163
$ws.Range("Q4").Value = 148.25351230029
$ws.Range("R4").Value = 1334.28161070261
$ws.Range("S4").Value = 0.09935246329884775
$ws.Range("T4").Value = 0.09935246329884777

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.730603666666667
$ws.Range("H5").Value = 11.191811
$ws.Range("I5").Value = 0.3405145117553424
$ws.Range("J5").Value = 0.3405145117553424
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.7065936666666666
$ws.Range("N5").Value = 2.119781
$ws.Range("O5").Value = 0.005187843618793344
$ws.Range("P5").Value = 0.005187843618793344
$ws.Range("Q5").Value = 2.636020923710111
$ws.Range("R5").Value = 23.724188313391
$ws.Range("S5").Value = 0.001766536036916484
$ws.Range("T5").Value = 0.001766536036916484

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.730603666666667
$ws.Range("H6").Value = 11.191811
$ws.Range("I6").Value = 0.3405145117553424
$ws.Range("J6").Value = 0.3405145117553424
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 111.9320066666667
$ws.Range("N6").Value = 335.79602
$ws.Range("O6").Value = 0.8218100075305903
$ws.Range("P6").Value = 0.8218100075305903
$ws.Range("Q6").Value = 417.5739544880245
$ws.Range("R6").Value = 3758.16559039222
$ws.Range("S6").Value = 0.2798382334699333
$ws.Range("T6").Value = 0.2798382334699333

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.730603666666667
$ws.Range("H7").Value = 11.191811
$ws.Range("I7").Value = 0.3405145117553424
$ws.Range("J7").Value = 0.3405145117553424
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 23.563205
$ws.Range("N7").Value = 70.689615
$ws.Range("O7").Value = 0.1730021488506163
$ws.Range("P7").Value = 0.1730021488506163
$ws.Range("Q7").Value = 87.90497897141833
$ws.Range("R7").Value = 791.144810742765
$ws.Range("S7").Value = 0.05890974224849268
$ws.Range("T7").Value = 0.05890974224849269

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.933443
$ws.Range("H8").Value = 2.800329
$ws.Range("I8").Value = 0.08520092612262004
$ws.Range("J8").Value = 0.08520092612262004
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.7065936666666666
$ws.Range("N8").Value = 2.119781
$ws.Range("O8").Value = 0.005187843618793344
$ws.Range("P8").Value = 0.005187843618793344
$ws.Range("Q8").Value = 0.6595649119943332
$ws.Range("R8").Value = 5.936084207948999
$ws.Range("S8").Value = 0.0004420090809005175
$ws.Range("T8").Value = 0.0004420090809005175

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.933443
$ws.Range("H9").Value = 2.800329
$ws.Range("I9").Value = 0.08520092612262004
$ws.Range("J9").Value = 0.08520092612262004
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 111.9320066666667
$ws.Range("N9").Value = 335.79602
$ws.Range("O9").Value = 0.8218100075305903
$ws.Range("P9").Value = 0.8218100075305903
$ws.Range("Q9").Value = 104.4821480989533
$ws.Range("R9").Value = 940.3393328905801
$ws.Range("S9").Value = 0.07001897373844364
$ws.Range("T9").Value = 0.07001897373844364

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.933443
$ws.Range("H10").Value = 2.800329
$ws.Range("I10").Value = 0.08520092612262004
$ws.Range("J10").Value = 0.08520092612262004
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.563205
$ws.Range("N10").Value = 70.689615
$ws.Range("O10").Value = 0.1730021488506163
$ws.Range("P10").Value = 0.1730021488506163
$ws.Range("Q10").Value = 21.994908764815
$ws.Range("R10").Value = 197.954178883335
$ws.Range("S10").Value = 0.01473994330327587
$ws.Range("T10").Value = 0.01473994330327588
